$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A3").Value = "words"

$ws2.Range("A4").Value = "id"
$ws2.Range("B4").Value = "lang1"
$ws2.Range("C4").Value = "lang2"
$ws2.Range("D4").Value = "sector"

$ws2.Range("A8").Value = "progress"

$ws2.Range("A9").Value = "user_id"
$ws2.Range("B9").Value = "words_id"
$ws2.Range("C9").Value = "percentage"

$ws2.Range("A16").Value = "дальше:"
$ws2.Range("A17").Value = "добавить кнопу +1 слово на изучение в шаблон game"
